$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F18").Value = "-"
$ws.Range("F19").Value = "-"
$ws.Range("F20").Value = "-"
$ws.Range("F21").Value = "-"
